$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2576.4
$ws.Range("I29").Value = 470.75
$ws.Range("J29").Value = 10999
$ws.Range("K29").Value = 1412.25
$ws.Range("L29").Value = 32997
$ws.Range("M29").Value = -1131.25
$ws.Range("N29").Value = -33559
$ws.Range("H38").Value = 118.416664
$ws.Range("I38").Value = 118.416664
$ws.Range("K38").Value = 355.249992
$ws.Range("M38").Value = 16.75000799999998
$ws.Range("H40").Value = 3970.625
$ws.Range("I40").Value = 3980
$ws.Range("J40").Value = 3969.2856
$ws.Range("K40").Value = 3980
$ws.Range("L40").Value = 3969.2856
$ws.Range("M40").Value = -3805
$ws.Range("N40").Value = -4319.2856
$ws.Range("H42").Value = 2106.3
$ws.Range("I42").Value = 370.8
$ws.Range("J42").Value = 3841.8
$ws.Range("K42").Value = 1112.4
$ws.Range("L42").Value = 11525.4
$ws.Range("M42").Value = -882.4000000000001
$ws.Range("N42").Value = -11985.4
$ws.Range("H43").Value = 3859.8
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 3859.8
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 3859.8
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -3997.8
$ws.Range("H62").Value = 14272.238
$ws.Range("I62").Value = 11095.637
$ws.Range("J62").Value = 17766.5
$ws.Range("K62").Value = 11095.637
$ws.Range("L62").Value = 17766.5
$ws.Range("M62").Value = -10471.637
$ws.Range("N62").Value = -19014.5
$ws.Range("H65").Value = 14272.238
$ws.Range("I65").Value = 11095.637
$ws.Range("J65").Value = 17766.5
$ws.Range("K65").Value = 55478.185
$ws.Range("L65").Value = 88832.5
$ws.Range("M65").Value = -52358.185
$ws.Range("N65").Value = -95072.5
$ws.Range("H74").Value = 5342.143
$ws.Range("I74").Value = 5342.143
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 5342.143
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -4406.143
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 5342.143
$ws.Range("I77").Value = 5342.143
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 26710.715
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -22030.715
$ws.Range("N77").ClearContents()
$ws.Range("H86").Value = 127230.25
$ws.Range("I86").Value = 145091.72
$ws.Range("K86").Value = 145091.72
$ws.Range("M86").Value = -143968.72
$ws.Range("H87").Value = 43633.332
$ws.Range("I87").Value = 20000
$ws.Range("J87").Value = 90900
$ws.Range("K87").Value = 20000
$ws.Range("L87").Value = 90900
$ws.Range("M87").Value = -18752
$ws.Range("N87").Value = -93396
$ws.Range("H89").Value = 127230.25
$ws.Range("I89").Value = 145091.72
$ws.Range("K89").Value = 725458.6
$ws.Range("M89").Value = -719842.6
$ws.Range("H90").Value = 43633.332
$ws.Range("I90").Value = 20000
$ws.Range("J90").Value = 90900
$ws.Range("K90").Value = 60000
$ws.Range("L90").Value = 272700
$ws.Range("M90").Value = -53760
$ws.Range("N90").Value = -285180
$ws.Range("H100").Value = 2678.8948
$ws.Range("I100").Value = 1799.4546
$ws.Range("K100").Value = 1799.4546
$ws.Range("M100").Value = -1258.4546
$ws.Range("H112").Value = 2865.4375
$ws.Range("I112").Value = 1261
$ws.Range("J112").Value = 3031.4138
$ws.Range("K112").Value = 3783
$ws.Range("L112").Value = 9094.241399999999
$ws.Range("M112").Value = -2675
$ws.Range("N112").Value = -11310.2414
$ws.Range("H113").Value = 8430.333000000001
$ws.Range("I113").Value = 11079.429
$ws.Range("J113").Value = 6112.375
$ws.Range("K113").Value = 11079.429
$ws.Range("L113").Value = 6112.375
$ws.Range("M113").Value = -7825.429
$ws.Range("N113").Value = -12620.375
$ws.Range("H138").Value = 2176.9333
$ws.Range("I138").Value = 1179.5358
$ws.Range("J138").Value = 3819.7058
$ws.Range("K138").Value = 3538.6074
$ws.Range("L138").Value = 11459.1174
$ws.Range("M138").Value = 1601.3926
$ws.Range("N138").Value = -21739.1174
$ws.Range("H141").Value = 1686.3721
$ws.Range("I141").Value = 1710.4524
$ws.Range("J141").Value = 675
$ws.Range("K141").Value = 5131.357199999999
$ws.Range("L141").Value = 2025
$ws.Range("M141").Value = 48.64280000000053
$ws.Range("N141").Value = -12385

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15387.821
$ws.Range("I32").Value = 14869.481
$ws.Range("K32").Value = 14869.481
$ws.Range("M32").Value = -14582.481
$ws.Range("H61").Value = 2545.5833
$ws.Range("I61").Value = 2085.7144
$ws.Range("J61").Value = 3189.4
$ws.Range("K61").Value = 2085.7144
$ws.Range("L61").Value = 3189.4
$ws.Range("M61").Value = -1873.7144
$ws.Range("N61").Value = -3613.4
$ws.Range("H132").Value = 25598.904
$ws.Range("I132").Value = 27954.395
$ws.Range("J132").Value = 3221.75
$ws.Range("K132").Value = 83863.185
$ws.Range("L132").Value = 9665.25
$ws.Range("M132").Value = -81333.185
$ws.Range("N132").Value = -14725.25
$ws.Range("H136").Value = 2545.5833
$ws.Range("I136").Value = 2085.7144
$ws.Range("J136").Value = 3189.4
$ws.Range("K136").Value = 6257.1432
$ws.Range("L136").Value = 9568.200000000001
$ws.Range("M136").Value = -3707.1432
$ws.Range("N136").Value = -14668.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 555.8
$ws.Range("I16").Value = 564.5263
$ws.Range("J16").Value = 390
$ws.Range("K16").Value = 564.5263
$ws.Range("L16").Value = 390
$ws.Range("M16").Value = -277.5263
$ws.Range("N16").Value = -964
$ws.Range("H55").Value = 14021.8
$ws.Range("I55").Value = 14021.8
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 14021.8
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -13706.8
$ws.Range("N55").ClearContents()
$ws.Range("H99").Value = 7834.4116
$ws.Range("I99").Value = 11356.7
$ws.Range("K99").Value = 11356.7
$ws.Range("M99").Value = -9858.700000000001
$ws.Range("H107").Value = 392.48
$ws.Range("I107").Value = 395.0909
$ws.Range("J107").Value = 373.33334
$ws.Range("K107").Value = 395.0909
$ws.Range("L107").Value = 373.33334
$ws.Range("M107").Value = 1524.9091
$ws.Range("N107").Value = -4213.33334
$ws.Range("H113").Value = 555.8
$ws.Range("I113").Value = 564.5263
$ws.Range("J113").Value = 390
$ws.Range("K113").Value = 564.5263
$ws.Range("L113").Value = 390
$ws.Range("M113").Value = 1605.4737
$ws.Range("N113").Value = -4730
$ws.Range("H122").Value = 1325.2307
$ws.Range("I122").Value = 1365.3636
$ws.Range("J122").Value = 1104.5
$ws.Range("K122").Value = 4096.0908
$ws.Range("L122").Value = 3313.5
$ws.Range("M122").Value = -1646.0908
$ws.Range("N122").Value = -8213.5
$ws.Range("H126").Value = 7834.4116
$ws.Range("I126").Value = 11356.7
$ws.Range("K126").Value = 34070.10000000001
$ws.Range("M126").Value = -31600.10000000001
$ws.Range("H134").Value = 20721.5
$ws.Range("I134").Value = 26504.732
$ws.Range("J134").Value = 2482.077
$ws.Range("K134").Value = 79514.196
$ws.Range("L134").Value = 7446.231000000001
$ws.Range("M134").Value = -76979.196
$ws.Range("N134").Value = -12516.231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 30.526316
$ws.Range("I2").Value = 36
$ws.Range("J2").Value = 18.666666
$ws.Range("K2").Value = 216
$ws.Range("L2").Value = 111.999996
$ws.Range("M2").Value = -103
$ws.Range("N2").Value = -337.999996
$ws.Range("H46").Value = 725.3570999999999
$ws.Range("I46").Value = 550.38464
$ws.Range("K46").Value = 1651.15392
$ws.Range("M46").Value = -1560.15392

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6002.8696
$ws.Range("I46").Value = 11116.333
$ws.Range("J46").Value = 2715.6428
$ws.Range("K46").Value = 11116.333
$ws.Range("L46").Value = 2715.6428
$ws.Range("M46").Value = -10928.333
$ws.Range("N46").Value = -3091.6428
$ws.Range("H55").Value = 1161.8636
$ws.Range("I55").Value = 815.9231
$ws.Range("J55").Value = 1661.5555
$ws.Range("K55").Value = 815.9231
$ws.Range("L55").Value = 1661.5555
$ws.Range("M55").Value = -642.9231
$ws.Range("N55").Value = -2007.5555
$ws.Range("H132").Value = 46270.605
$ws.Range("I132").Value = 56995
$ws.Range("K132").Value = 170985
$ws.Range("M132").Value = -168455
$ws.Range("H133").Value = 63999
$ws.Range("J133").Value = 63999
$ws.Range("L133").Value = 63999
$ws.Range("N133").Value = -69059

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 54999
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H64").Value = 58100
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 58100
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H136").Value = 1669.0435
$ws.Range("I136").Value = 1599.8823
$ws.Range("J136").Value = 1865
$ws.Range("K136").Value = 4799.6469
$ws.Range("L136").Value = 5595
$ws.Range("M136").Value = -2249.6469
$ws.Range("N136").Value = -10695
